# Actualización automática del mapa: agrega el nuevo caso -567 como fila 76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    # Forzar texto para que los valores numéricos (OT, Comuna, Caso, fecha, etc.)
    # se guarden como cadena y no se interpreten como número/fecha.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    # Quitar el formato aplicado para que la celda quede sin estilo, igual
    # que el resto de las filas de datos de la hoja.
    $cell.Style = "Normal"
}

$newRow = 76

Set-TextCell $newRow 1  "-567"
Set-TextCell $newRow 2  "8/25/2025"
Set-TextCell $newRow 3  "Franco 2543"
Set-TextCell $newRow 4  "12"
Set-TextCell $newRow 5  "809184735"
Set-TextCell $newRow 6  "NEW"
Set-TextCell $newRow 7  "Pendiente"
Set-TextCell $newRow 8  "Cambio"
$ws.Cells.Item($newRow, 9).Value = 0
Set-TextCell $newRow 10 "Cambio"
Set-TextCell $newRow 11 "Sin equipos"
Set-TextCell $newRow 12 "Pasante"
$ws.Cells.Item($newRow, 13).Value = -58.502565
$ws.Cells.Item($newRow, 14).Value = -34.578977
Set-TextCell $newRow 15 "Paternal"
Set-TextCell $newRow 16 "Capital Norte"
